$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("M2").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("N2").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M3").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("M4").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N4").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M5").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N6").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M7").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("M8").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N8").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M9").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N9").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M10").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("M11").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("M12").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N12").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M13").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("N13").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M14").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("M15").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N15").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M16").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("M17").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N17").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("N18").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M19").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("N19").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M20").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N20").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M21").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N21").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("N22").Value2 = "Sí, el uso de redes sociales ha tenido un impacto negativo en mi experiencia académica"
$ws.Range("M23").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("M24").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N24").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M25").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("N26").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M27").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("N27").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M28").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("M29").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("N29").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M30").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N30").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M31").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N31").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M32").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N32").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M33").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N33").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M34").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("M35").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N36").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M37").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("M38").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("N38").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M39").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N39").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("N40").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M41").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N41").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("N43").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M44").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("M45").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N45").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M47").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N47").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M48").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N48").Value2 = "Sí, el uso de redes sociales ha tenido un impacto negativo en mi experiencia académica"
$ws.Range("N49").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M51").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N51").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M52").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("M53").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N53").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M54").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N54").Value2 = "Sí, el uso de redes sociales ha tenido un impacto negativo en mi experiencia académica"
$ws.Range("M56").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N56").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("N57").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M58").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("M59").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N59").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M60").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N60").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M61").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N61").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M62").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N63").Value2 = "Sí, el uso de redes sociales ha tenido un impacto negativo en mi experiencia académica"
$ws.Range("M64").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("M65").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("M67").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N67").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M68").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N68").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M69").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N70").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M71").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("N71").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M72").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N72").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M73").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N74").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M75").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N75").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M76").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N76").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M77").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N77").Value2 = "Sí, el uso de redes sociales ha tenido un impacto negativo en mi experiencia académica"
$ws.Range("M78").Value2 = "No, creo que puedo controlar mi uso de redes sociales y no afecta mi vida académica"
$ws.Range("M79").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N79").Value2 = "Sí, el uso de redes sociales ha tenido un impacto positivo en mi experiencia académica"
$ws.Range("M80").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("M81").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N81").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M84").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N86").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M87").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("M89").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("M90").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("N90").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M92").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("M94").Value2 = "No uso redes sociales o las uso muy poco"
$ws.Range("N95").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("N96").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"
$ws.Range("M98").Value2 = "Sí, creo que paso demasiado tiempo en redes sociales y me afecta negativamente"
$ws.Range("M99").Value2 = "Sí, a veces me cuesta limitar mi tiempo en redes sociales cuando debería estar haciendo otras cosas"
$ws.Range("N100").Value2 = "No estoy seguro/a del impacto del uso de redes sociales en mi experiencia académica"

$ws.Range("P8").Select()

Write-Host "Applied 128 cell updates"
